$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns right by one.
$ws.Range("A:A").Insert()

# Determine the last used row (data rows), based on the sheet's used range.
$lastRow = $ws.Cells.SpecialCells(11).Row

# Copy the header cell's format (bold, centered, bordered) so it can be applied
# to the new index column cells.
$ws.Range("B1").Copy()

# Fill the new column A with a zero-based sequential index for each data row
# (rows 2..lastRow), using the same style as the header row cells.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.PasteSpecial(-4122)
    $cell.Value = $r - 2
}

$excel.CutCopyMode = 0
